$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C4 from "okay" to "Done"
$ws.Range("C4").Value = "Done"

# Remove the extra rows (5-8) that contained leftover/test data
$ws.Range("A5:C8").EntireRow.Delete()
